# Updated cryptos list - price/volume refresh + Litecoin/Uniswap rank swap (rows 20-21)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.916.19"
$ws.Range("E2").Value = "  -0.44%  "

# Row 3
$ws.Range("D3").Value = "2.235.41"
$ws.Range("E3").Value = "  -0.59%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "'272.57"
$ws.Range("E5").Value = "  +5.55%  "

# Row 6
$ws.Range("D6").Value = "'88.13"
$ws.Range("E6").Value = "  +8.49%  "

# Row 7
$ws.Range("D7").Value = "'0.622"
$ws.Range("E7").Value = "  -0.05%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("D9").Value = "'0.604"
$ws.Range("E9").Value = "  -0.34%  "

# Row 10
$ws.Range("D10").Value = "'45.06"
$ws.Range("E10").Value = "  +2.82%  "

# Row 11
$ws.Range("D11").Value = "'0.0919"
$ws.Range("E11").Value = "  -1.65%  "

# Row 12
$ws.Range("D12").Value = "'7.65"
$ws.Range("E12").Value = "  +7.83%  "

# Row 13
$ws.Range("D13").Value = "'0.104"
$ws.Range("E13").Value = "  +0.67%  "

# Row 14
$ws.Range("D14").Value = "2.565.82"
$ws.Range("E14").Value = "  -0.54%  "

# Row 15
$ws.Range("D15").Value = "'15.08"
$ws.Range("E15").Value = "  +2.07%  "

# Row 16
$ws.Range("D16").Value = "2.227.77"
$ws.Range("E16").Value = "  -1.83%  "

# Row 17
$ws.Range("D17").Value = "'0.795"
$ws.Range("E17").Value = "  -0.37%  "

# Row 18
$ws.Range("D18").Value = "43.832.51"
$ws.Range("E18").Value = "  -0.46%  "

# Row 19
$ws.Range("D19").Value = "'0.0000104"
$ws.Range("E19").Value = "  -1.26%  "

# Row 20
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'5.97"
$ws.Range("E20").Value = "  -1.80%  "

# Row 21
$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").Value = "'70.17"
$ws.Range("E21").Value = "  -1.96%  "

# Row 22
$ws.Range("D22").Value = "'2.34"
$ws.Range("E22").Value = "  +0.29%  "

# Row 23
$ws.Range("D23").Value = "'232.43"
$ws.Range("E23").Value = "  -0.38%  "

# Row 24
$ws.Range("D24").Value = "'8.68"
$ws.Range("E24").Value = "  -7.72%  "

# Row 25
$ws.Range("E25").Value = "  -0.07%  "

# Row 26
$ws.Range("D26").Value = "'2.54"
$ws.Range("E26").Value = "  +13.68%  "

# Row 27
$ws.Range("D27").Value = "'10.83"
$ws.Range("E27").Value = "  -0.67%  "

# Row 28
$ws.Range("D28").Value = "'3.52"
$ws.Range("E28").Value = "  +4.22%  "

# Row 29
$ws.Range("E29").Value = "  +3.67%  "

# Row 30
$ws.Range("D30").Value = "'38.94"
$ws.Range("E30").Value = "  -3.60%  "

# Row 31
$ws.Range("D31").Value = "'172.87"
$ws.Range("E31").Value = "  -0.08%  "

# Row 32
$ws.Range("D32").Value = "'0.0913"
$ws.Range("E32").Value = "  +2.86%  "

# Row 33
$ws.Range("D33").Value = "'20.76"
$ws.Range("E33").Value = "  +0.38%  "

# Row 34
$ws.Range("D34").Value = "'5.33"
$ws.Range("E34").Value = "  +0.13%  "

# Row 35
$ws.Range("E35").Value = "  +0.01%  "

# Row 36
$ws.Range("D36").Value = "'0.111"
$ws.Range("E36").Value = "  -1.73%  "

# Row 37
$ws.Range("E37").Value = "  -3.44%  "

# Row 38
$ws.Range("D38").Value = "'4.26"
$ws.Range("E38").Value = "  -6.49%  "

# Row 39
$ws.Range("E39").Value = "  +17.11%  "

# Row 40
$ws.Range("D40").Value = "'2.17"
$ws.Range("E40").Value = "  +1.15%  "

# Row 41
$ws.Range("D41").Value = "'12.42"
$ws.Range("E41").Value = "  -3.73%  "

# Row 42
$ws.Range("D42").Value = "'0.213"
$ws.Range("E42").Value = "  +4.13%  "

# Row 43
$ws.Range("D43").Value = "'63.44"
$ws.Range("E43").Value = "  -0.32%  "

# Row 44
$ws.Range("D44").Value = "'5.39"
$ws.Range("E44").Value = "  -3.38%  "

# Row 45
$ws.Range("E45").Value = "  -0.64%  "

# Row 46
$ws.Range("D46").Value = "'0.0984"
$ws.Range("E46").Value = "  -0.50%  "

# Row 47
$ws.Range("D47").Value = "'99.97"
$ws.Range("E47").Value = "  -4.25%  "

# Row 48
$ws.Range("E48").Value = "  +0.53%  "

# Row 49
$ws.Range("E49").Value = "  +3.04%  "

# Row 50
$ws.Range("E50").Value = "  -4.53%  "

# Row 51
$ws.Range("D51").Value = "'1.49"
$ws.Range("E51").Value = "  -2.29%  "
